$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "14-jul" date column header in AD1
$ws.Range("AD1").Value = "14-jul"

# Add the AD values for rows 2-18 (mirrors the AC column's daily data series)
$adValues = @(
    0,
    11.826604548358478,
    15.3530921747726,
    28.248194983727391,
    0,
    3.0600308335208122,
    12.92637578567901,
    22.263885348326731,
    20.379450246782088,
    10.958627514887576,
    0,
    9.9777776693361329,
    0,
    0,
    15.45784680051514,
    0,
    0
)

for ($i = 0; $i -lt $adValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 30).Value = $adValues[$i]
}

# Move the active selection to AE8 (matches the recorded sheet view state)
$ws.Range("AE8").Select()
